## "Se actualizo la lista de paquetes"
## Refresh the package list table (columns A=BASE, B=MODELADO, C=SECUNDARIOS).
## Column C (SECUNDARIOS) is unchanged by this edit; only A and B get new
## contents, and the table shrinks by one row (18 -> 17 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BASE column (A2:A17)
$colA = @(
    "corrplot", "ggthemes", "knitr", "openxlsx", "DataExplorer", "class",
    "visdat", "ggvis", "GGally", "fastDummies", "gridExtra", "ROCR",
    "pROC", "reshape2", "devtools", "tidyverse"
)
for ($i = 0; $i -lt $colA.Count; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $colA[$i]
}

# New MODELADO column (B2:B9); B10 and B11 no longer have a value.
$colB = @("caret", "Metrics", "nortest", "lmtest", "glmnet", "klaR", "car", "tensorflow")
for ($i = 0; $i -lt $colB.Count; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}
$ws.Cells.Item(10, 2).Value = $null
$ws.Cells.Item(11, 2).Value = $null

# Column C (SECUNDARIOS) is untouched -- leave as-is.

# The table lost its last row (was 18 rows, now 17) -- remove it outright so
# no stray blank row / stale dimension is left behind.
$ws.Rows(18).Delete()

# Keep the underline styling on "reshape2" (now row 15).
$ws.Cells.Item(15, 1).Font.Underline = $true

# Match the recorded selection from the edit.
$ws.Range("B15").Select()
